# Generate Report for handoff
# Updates the "Latest Handoff Datetime" (column D) for the 4th data row
# (the 66958070-8c4b-4e5a-ae0e-aef9d9d3421a file) on both the "zh-cn" and
# "de-de" localization-status worksheets, reflecting a newer handoff run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D4").Value = "2016-01-18 02:41:35"
$wsDeDe.Range("D4").Value = "2016-01-18 02:41:47"
